{"js": "// Peque\u00f1as correcciones modelo verbal\n// 1) \"dependencia a la que pertenece\" -> \"departamento al que pertenece\"\n// 2) \"Dependencia, nombre del proyecto...\" -> \"Departamento, nombre del proyecto...\"\n// 3) Move the \"_GoBack\" bookmark from before \"disponibilidad horaria...\" to\n//    sit inside the newly-typed word \"Departamento\" (right after \"Departa\"),\n//    matching where Word leaves the last-edit marker.\n\nconst body = context.document.body;\n\n// --- Change 1: \"dependencia a la \" -> \"departamento al \" ---------------\nlet results = body.search(\"dependencia a la \", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for 'dependencia a la ', found \" + results.items.length);\n}\nresults.items[0].insertText(\"departamento al \", \"Replace\");\nawait context.sync();\n\n// --- Change 2: \"Dependencia, nombre del proyecto\" -> \"Departamento, nombre del proyecto\" ---\nresults = body.search(\"Dependencia, nombre del proyecto\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for 'Dependencia, nombre del proyecto', found \" + results.items.length);\n}\nresults.items[0].insertText(\"Departamento, nombre del proyecto\", \"Replace\");\nawait context.sync();\n\n// --- Change 3: relocate the \"_GoBack\" bookmark --------------------------\n// Remove the existing \"_GoBack\" bookmark (currently sitting right before\n// \"disponibilidad horaria, porcentaje de avance...\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-insert it collapsed right after \"Departa\" (i.e. inside \"Departamento\"),\n// which is where Word records the point of the last text edit.\nresults = body.search(\"Departa\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for 'Departa', found \" + results.items.length);\n}\nconst insertionPoint = results.items[0].getRange(\"End\");\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Peque\u00f1as correcciones modelo verbal\n# 1) \"dependencia a la que pertenece\" -> \"departamento al que pertenece\"\n# 2) \"Dependencia, nombre del proyecto...\" -> \"Departamento, nombre del proyecto...\"\n# 3) Move the \"_GoBack\" bookmark from before \"disponibilidad horaria...\" to\n#    sit inside the newly-typed word \"Departamento\" (right after \"Departa\"),\n#    matching where Word leaves the last-edit marker.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: \"dependencia a la \" -> \"departamento al \" ---------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"dependencia a la \"\n$find1.Replacement.Text = \"departamento al \"\n$find1.MatchCase = $true\n$find1.MatchWholeWord = $false\n$find1.Forward = $true\n$find1.Wrap = 1  # wdFindContinue\n$find1.Execute($find1.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# --- Change 2: \"Dependencia, nombre del proyecto\" -> \"Departamento, nombre del proyecto\" ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Dependencia, nombre del proyecto\"\n$find2.Replacement.Text = \"Departamento, nombre del proyecto\"\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Execute($find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# --- Change 3: relocate the \"_GoBack\" bookmark --------------------------\n# Remove the existing \"_GoBack\" bookmark (currently sitting right before\n# \"disponibilidad horaria, porcentaje de avance...\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Re-insert it collapsed right after \"Departa\" (i.e. inside \"Departamento\"),\n# which is where Word records the point of the last text edit.\n$searchRange = $d.Content\n$find3 = $searchRange.Find\n$find3.ClearFormatting()\n$find3.Text = \"Departa\"\n$find3.MatchCase = $true\n$find3.Forward = $true\n$find3.Wrap = 1\n$found = $find3.Execute()\nif (-not $found) {\n  throw \"Could not find 'Departa' to relocate the _GoBack bookmark\"\n}\n$searchRange.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $searchRange) | Out-Null\n"}
